$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg).
# These correspond to a reshuffle of the existing rows' data.
$rows = @{
    2  = @(44762, 50, 2300, 2300, 2300, 2300)
    3  = @(44473, 120, 1200, 1200, 1200, 1200)
    4  = @(44431, 100, 1300, 1300, 1300, 1300)
    5  = @(44753, 160, 2300, 2300, 2300, 2300)
    6  = @(44749, 120, 2300, 2300, 2300, 2300)
    7  = @(44357, 35, 1000, 1000, 1000, 1000)
    8  = @(44812, 50, 2500, 2500, 2500, 2500)
    9  = @(44811, 60, 2500, 2500, 2500, 2500)
    10 = @(44435, 130, 1300, 1300, 1300, 1300)
    11 = @(44343, 60, 1300, 1300, 1300, 1300)
    12 = @(44424, 50, 1200, 1200, 1200, 1200)
    13 = @(44760, 80, 2300, 2300, 2300, 2300)
    14 = @(44476, 80, 1200, 1200, 1200, 1200)
    15 = @(44418, 40, 1200, 1200, 1200, 1200)
    16 = @(44748, 300, 2300, 2300, 2300, 2300)
    17 = @(44432, 30, 1300, 1300, 1300, 1300)
    18 = @(44830, 50, 2500, 2500, 2500, 2500)
    19 = @(44405, 50, 1200, 1200, 1200, 1200)
    20 = @(44417, 80, 1200, 1200, 1200, 1200)
    21 = @(44763, 50, 2300, 2300, 2300, 2300)
    22 = @(44438, 60, 1200, 1200, 1200, 1200)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("M$r").Value = $vals[1]
    $ws.Range("N$r").Value = $vals[2]
    $ws.Range("O$r").Value = $vals[3]
    $ws.Range("P$r").Value = $vals[4]
    $ws.Range("S$r").Value = $vals[5]
}
